$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "32×27=" "72×18="
Replace-Text "67×20=" "81×33="
Replace-Text "55×40=" "25×34="
Replace-Text "44×59=" "50×43="
Replace-Text "46×57=" "24×21="
Replace-Text "25×79=" "60×51="
Replace-Text "56×17=" "46×38="
Replace-Text "91×63=" "54×87="
Replace-Text "44×75=" "22×18="
Replace-Text "95×60=" "52×32="
Replace-Text "68×28=" "35×25="
Replace-Text "85×14=" "58×29="
Replace-Text "64×49=" "72×48="
Replace-Text "82×24=" "17×55="
Replace-Text "27×66=" "73×78="
Replace-Text "98×24=" "74×25="
Replace-Text "47×76=" "67×80="
Replace-Text "63×21=" "69×54="
Replace-Text "47×94=" "89×70="
Replace-Text "35×76=" "39×91="
Replace-Text "13×13=" "75×20="
Replace-Text "30×19=" "18×35="
Replace-Text "73×29=" "31×54="
Replace-Text "25×24=" "23×23="
Replace-Text "39×51=" "29×27="
